$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content edit: "offcanvas" -> "off canvas" (cell A24)
$ws.Range("A24").Value = "off canvas"

# View edit: selection moved to D6 (and sheet scrolled back to top-left)
[void]$ws.Range("D6").Select()
